# Applies the "Updated symbol list" data refresh (Fri Jan 20 07:07:13 UTC 2023)
# to the cryptos sheet: refreshed Price (D) / Volume 1h (E) values and bumped
# the "Hora" (G) column from 6 -> 7 for every data row (2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = '288.10' },
    @{ Cell = "E2"; Value = '-0.86%' },
    @{ Cell = "G2"; Value = '7' },
    @{ Cell = "D3"; Value = '31.03' },
    @{ Cell = "E3"; Value = '1.00%' },
    @{ Cell = "G3"; Value = '7' },
    @{ Cell = "D4"; Value = '4.932' },
    @{ Cell = "E4"; Value = '-0.51%' },
    @{ Cell = "G4"; Value = '7' },
    @{ Cell = "D5"; Value = '0.07364' },
    @{ Cell = "E5"; Value = '1.99%' },
    @{ Cell = "G5"; Value = '7' },
    @{ Cell = "D6"; Value = '2.260' },
    @{ Cell = "E6"; Value = '25.70%' },
    @{ Cell = "G6"; Value = '7' },
    @{ Cell = "D7"; Value = '7.733' },
    @{ Cell = "E7"; Value = '0.77%' },
    @{ Cell = "G7"; Value = '7' },
    @{ Cell = "D8"; Value = '3.736' },
    @{ Cell = "E8"; Value = '0.24%' },
    @{ Cell = "G8"; Value = '7' },
    @{ Cell = "D9"; Value = '0.9065' },
    @{ Cell = "E9"; Value = '1.14%' },
    @{ Cell = "G9"; Value = '7' },
    @{ Cell = "D10"; Value = '0.08696' },
    @{ Cell = "E10"; Value = '13.60%' },
    @{ Cell = "G10"; Value = '7' },
    @{ Cell = "D11"; Value = '0.1683' },
    @{ Cell = "E11"; Value = '1.87%' },
    @{ Cell = "G11"; Value = '7' },
    @{ Cell = "D12"; Value = '0.08226' },
    @{ Cell = "E12"; Value = '2.60%' },
    @{ Cell = "G12"; Value = '7' },
    @{ Cell = "D13"; Value = '0.03115' },
    @{ Cell = "E13"; Value = '2.62%' },
    @{ Cell = "G13"; Value = '7' },
    @{ Cell = "D14"; Value = '0.09933' },
    @{ Cell = "E14"; Value = '-0.89%' },
    @{ Cell = "G14"; Value = '7' },
    @{ Cell = "D15"; Value = '0.001503' },
    @{ Cell = "E15"; Value = '-0.18%' },
    @{ Cell = "G15"; Value = '7' },
    @{ Cell = "D16"; Value = '0.005761' },
    @{ Cell = "E16"; Value = '-0.24%' },
    @{ Cell = "G16"; Value = '7' },
    @{ Cell = "D17"; Value = '3.490' },
    @{ Cell = "E17"; Value = '0.81%' },
    @{ Cell = "G17"; Value = '7' },
    @{ Cell = "D18"; Value = '2.097' },
    @{ Cell = "E18"; Value = '0.79%' },
    @{ Cell = "G18"; Value = '7' },
    @{ Cell = "D19"; Value = '0.3326' },
    @{ Cell = "E19"; Value = '0.33%' },
    @{ Cell = "G19"; Value = '7' },
    @{ Cell = "D20"; Value = '0.1294' },
    @{ Cell = "E20"; Value = '-1.59%' },
    @{ Cell = "G20"; Value = '7' },
    @{ Cell = "E21"; Value = '-5.14%' },
    @{ Cell = "G21"; Value = '7' },
    @{ Cell = "D22"; Value = '0.2124' },
    @{ Cell = "E22"; Value = '1.09%' },
    @{ Cell = "G22"; Value = '7' },
    @{ Cell = "D23"; Value = '0.04557' },
    @{ Cell = "E23"; Value = '0.84%' },
    @{ Cell = "G23"; Value = '7' },
    @{ Cell = "D24"; Value = '0.001209' },
    @{ Cell = "E24"; Value = '-0.51%' },
    @{ Cell = "G24"; Value = '7' },
    @{ Cell = "D25"; Value = '0.004140' },
    @{ Cell = "E25"; Value = '3.15%' },
    @{ Cell = "G25"; Value = '7' },
    @{ Cell = "E26"; Value = '4.05%' },
    @{ Cell = "G26"; Value = '7' },
    @{ Cell = "D27"; Value = '0.0003400' },
    @{ Cell = "E27"; Value = '-95.48%' },
    @{ Cell = "G27"; Value = '7' },
    @{ Cell = "G28"; Value = '7' },
    @{ Cell = "G29"; Value = '7' },
    @{ Cell = "G30"; Value = '7' },
    @{ Cell = "G31"; Value = '7' },
    @{ Cell = "G32"; Value = '7' },
    @{ Cell = "G33"; Value = '7' },
    @{ Cell = "G34"; Value = '7' },
    @{ Cell = "G35"; Value = '7' },
    @{ Cell = "G36"; Value = '7' },
    @{ Cell = "G37"; Value = '7' },
    @{ Cell = "G38"; Value = '7' },
    @{ Cell = "D39"; Value = '0.01576' },
    @{ Cell = "E39"; Value = '-1.68%' },
    @{ Cell = "G39"; Value = '7' },
    @{ Cell = "D40"; Value = '0.04464' },
    @{ Cell = "E40"; Value = '1.65%' },
    @{ Cell = "G40"; Value = '7' },
    @{ Cell = "D41"; Value = '0.007349' },
    @{ Cell = "E41"; Value = '-0.07%' },
    @{ Cell = "G41"; Value = '7' },
    @{ Cell = "D42"; Value = '0.009580' },
    @{ Cell = "E42"; Value = '24.34%' },
    @{ Cell = "G42"; Value = '7' },
    @{ Cell = "D43"; Value = '0.1322' },
    @{ Cell = "E43"; Value = '1.04%' },
    @{ Cell = "G43"; Value = '7' },
    @{ Cell = "D44"; Value = '0.002242' },
    @{ Cell = "E44"; Value = '9.32%' },
    @{ Cell = "G44"; Value = '7' },
    @{ Cell = "D45"; Value = '0.008438' },
    @{ Cell = "E45"; Value = '-8.43%' },
    @{ Cell = "G45"; Value = '7' },
    @{ Cell = "D46"; Value = '0.00006115' },
    @{ Cell = "E46"; Value = '4.25%' },
    @{ Cell = "G46"; Value = '7' },
    @{ Cell = "E47"; Value = '0.12%' },
    @{ Cell = "G47"; Value = '7' },
    @{ Cell = "D48"; Value = '2.187' },
    @{ Cell = "E48"; Value = '-2.59%' },
    @{ Cell = "G48"; Value = '7' },
    @{ Cell = "D49"; Value = '0.002003' },
    @{ Cell = "E49"; Value = '-33.27%' },
    @{ Cell = "G49"; Value = '7' },
    @{ Cell = "D50"; Value = '0.00002103' },
    @{ Cell = "E50"; Value = '0.12%' },
    @{ Cell = "G50"; Value = '7' },
    @{ Cell = "D51"; Value = '0.0002003' },
    @{ Cell = "E51"; Value = '0.12%' },
    @{ Cell = "G51"; Value = '7' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text storage so numeric-looking strings (e.g. "288.10", "0.0003400")
    # keep their exact literal formatting instead of being coerced to a number.
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
